$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "egrggerge"
$ws.Range("A5").Value = "ergerghog24"
$ws.Range("A6").Value = "wrlrgpirhgwp"
$ws.Range("A8").Value = "oinohrg3wg"

$ws.Range("A8").Select()
